$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.158.97"
$ws.Range("E2").Value = "  +2.67%  "

$ws.Range("D3").Value = "3.563.22"
$ws.Range("E3").Value = "  +7.03%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.77"
$ws.Range("E5").Value = "  +3.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "637.90"
$ws.Range("E6").Value = "  +3.26%  "

$ws.Range("E7").Value = "  +7.60%  "

$ws.Range("E8").Value = "  +3.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  +10.27%  "

$ws.Range("D11").Value = "3.561.31"
$ws.Range("E11").Value = "  +6.99%  "

$ws.Range("E12").Value = "  +2.79%  "

$ws.Range("E13").Value = "  +4.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.46"
$ws.Range("E14").Value = "  +8.50%  "

$ws.Range("D15").Value = "4.240.18"
$ws.Range("E15").Value = "  +7.44%  "

$ws.Range("D16").Value = "96.029.43"
$ws.Range("E16").Value = "  +2.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000254"
$ws.Range("E17").Value = "  +4.35%  "

$ws.Range("D18").Value = "3.563.07"
$ws.Range("E18").Value = "  +7.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.13"
$ws.Range("E19").Value = "  +20.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.93"
$ws.Range("E20").Value = "  -1.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.08"
$ws.Range("E21").Value = "  +5.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.503"
$ws.Range("E22").Value = "  +12.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "516.66"
$ws.Range("E23").Value = "  +4.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.43"
$ws.Range("E24").Value = "  -0.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000195"
$ws.Range("E25").Value = "  +7.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.70"
$ws.Range("E26").Value = "  +8.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.86"
$ws.Range("E27").Value = "  +5.97%  "

$ws.Range("E28").Value = "  +5.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.12"
$ws.Range("E29").Value = "  +19.60%  "

$ws.Range("E30").Value = "  +4.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.57"
$ws.Range("E31").Value = "  +4.96%  "

$ws.Range("E33").Value = "  +5.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.01"
$ws.Range("E34").Value = "  +1.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "30.22"
$ws.Range("E35").Value = "  +6.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.565"
$ws.Range("E36").Value = "  +6.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "579.61"
$ws.Range("E37").Value = "  +8.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.89"
$ws.Range("E38").Value = "  +6.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.50"
$ws.Range("E39").Value = "  +10.35%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.152"
$ws.Range("E40").Value = "  +2.93%  "

$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("E42").Value = "  +7.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0434"
$ws.Range("E43").Value = "  +4.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.76"
$ws.Range("E44").Value = "  +4.68%  "

$ws.Range("E45").Value = "  -0.83%  "

$ws.Range("E46").Value = "  +4.53%  "

$ws.Range("E47").Value = "  -3.12%  "

$ws.Range("E48").Value = "  +3.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.95"
$ws.Range("E49").Value = "  +3.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.19"
$ws.Range("E50").Value = "  +3.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.12"
$ws.Range("E51").Value = "  +2.79%  "
